# Add 2022-Q1 data:
#  - the existing "总计" sheet becomes "2022-Q1" and is refilled with the
#    quarter's fund-holdings detail (same shape as the other quarter sheets)
#  - a brand-new "总计" sheet is appended right after it, carrying the old
#    summary rows plus a new 2022-Q1 row on top
#
# NOTE: sheet handles returned by Worksheets.Add()/.Item() become stale for
# Cells addressing once another structural operation (Add/Move/rename of a
# *different* sheet) runs afterwards, so every structural change happens
# first and sheet references used for writing cell data are (re)fetched by
# name only after the sheet layout is final.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Structural changes only: rename + add + reposition
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

$q1ForMove = $wb.Worksheets.Item("2022-Q1")
$newTotalForMove = $wb.Worksheets.Item("总计")
$newTotalForMove.Move($null, $q1ForMove)

# ---------------------------------------------------------------------
# 2) Re-fetch final sheet references (layout is now settled) and a
#    "style 2" source cell to clone the bold header / index-column look
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")
$styleSource = $wb.Worksheets.Item("2021-Q4").Range("A2")

# ---------------------------------------------------------------------
# 3) "2022-Q1" sheet: fund holdings detail
# ---------------------------------------------------------------------
$q1.Cells.Clear()

$headers1 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers1) {
    $q1.Cells.Item(1, $col).Value = $h
    $styleSource.Copy()
    $q1.Cells.Item(1, $col).PasteSpecial(-4122)
    $col = $col + 1
}

$funds = @(
    @("007872", "金信稳健策略灵活配置混合", "25.57", "93.73", "7.55", "1.9305", 5),
    @("004666", "长城久嘉创新成长灵活配置混合", "26.19", "92.18", "7.37", "1.9302", 5),
    @("002256", "金信行业优选灵活配置混合", "2.43", "93.89", "7.81", "0.1898", 5),
    @("005434", "鹏华睿投灵活配置混合", "3.41", "82.48", "1.95", "0.0665", 8),
    @("005117", "金信价值精选灵活配置混合A", "0.83", "86.42", "5.91", "0.0491", 2),
    @("005296", "南华丰淳混合A", "1.04", "94.59", "3.72", "0.0387", 7),
    @("004223", "金信多策略精选灵活配置混合", "0.36", "93.14", "7.28", "0.0262", 3),
    @("006692", "金信消费升级股票A", "0.62", "94.05", "4.08", "0.0253", 7),
    @("004926", "中航军民融合精选混合A", "0.35", "91.27", "6.68", "0.0234", 2),
    @("005297", "南华丰淳混合C", "0.53", "94.59", "3.72", "0.0197", 7),
    @("004927", "中航军民融合精选混合C", "0.27", "91.27", "6.68", "0.0180", 2),
    @("009317", "金信核心竞争力灵活配置混合", "0.19", "89.48", "9.16", "0.0174", 4),
    @("005000", "泰康泉林量化价值精选混合A", "0.64", "93.61", "2.12", "0.0136", 7),
    @("002630", "江信瑞福灵活配置混合A", "0.52", "43.17", "2.50", "0.0130", 6),
    @("002631", "江信瑞福灵活配置混合C", "0.50", "43.17", "2.50", "0.0125", 6),
    @("002862", "金信量化精选灵活配置混合", "0.16", "94.28", "5.36", "0.0086", 5),
    @("002810", "金信转型创新成长灵活配置混合", "0.18", "81.12", "4.59", "0.0083", 3),
    @("006693", "金信消费升级股票C", "0.20", "94.05", "4.08", "0.0082", 7),
    @("005111", "泰康泉林量化价值精选混合C", "0.25", "93.61", "2.12", "0.0053", 7),
    @("005118", "金信价值精选灵活配置混合C", "0.05", "86.42", "5.91", "0.0030", 2)
)

$row = 2
$idx = 0
foreach ($f in $funds) {
    $q1.Cells.Item($row, 1).Value = $idx
    $styleSource.Copy()
    $q1.Cells.Item($row, 1).PasteSpecial(-4122)

    $q1.Cells.Item($row, 2).Value = "'" + $f[0]
    $q1.Cells.Item($row, 3).Value = $f[1]
    $q1.Cells.Item($row, 4).Value = "'" + $f[2]
    $q1.Cells.Item($row, 5).Value = "'" + $f[3]
    $q1.Cells.Item($row, 6).Value = "'" + $f[4]
    $q1.Cells.Item($row, 7).Value = "'" + $f[5]
    $q1.Cells.Item($row, 8).Value = $f[6]

    $row = $row + 1
    $idx = $idx + 1
}

# ---------------------------------------------------------------------
# 4) "总计" sheet: quarter-by-quarter summary (2022-Q1 on top)
# ---------------------------------------------------------------------
$headers2 = @("日期", "持有数量(只)", "持有市值(亿元)")
$col = 2
foreach ($h in $headers2) {
    $total.Cells.Item(1, $col).Value = $h
    $styleSource.Copy()
    $total.Cells.Item(1, $col).PasteSpecial(-4122)
    $col = $col + 1
}

$summary = @(
    @("2022-Q1", 20, 4.41),
    @("2021-Q4", 19, 13.16),
    @("2021-Q3", 26, 11.79),
    @("2021-Q2", 9, 2.35),
    @("2021-Q1", 12, 2.87),
    @("2020-Q4", 3, 0.82)
)

$row = 2
$idx = 0
foreach ($s in $summary) {
    $total.Cells.Item($row, 1).Value = $idx
    $styleSource.Copy()
    $total.Cells.Item($row, 1).PasteSpecial(-4122)

    $total.Cells.Item($row, 2).Value = $s[0]
    $total.Cells.Item($row, 3).Value = $s[1]
    $total.Cells.Item($row, 4).Value = $s[2]

    $row = $row + 1
    $idx = $idx + 1
}
